$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto price/volume refresh (GitHub Actions scheduled update)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.078.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.509"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.671.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.078.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +4.97%  "
$ws.Range("E23").Value = "  +4.49%  "
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.529.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("E35").Value = "  +8.69%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.580"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.891"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.94%  "
$ws.Range("E43").Value = "  +3.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.797.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.777"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("E51").Value = "  +3.27%  "
